$wb = $excel.ActiveWorkbook

# --- Fix existing sheet selections / tab state -------------------------------------------------
# Germany sheet (sheet1.xml): selection currently spans whole sheet, should be A1:D10
$wsGermany = $wb.Worksheets.Item("Germany")
$wsGermany.Range("A1:D10").Select() | Out-Null

# Swiss sheet (sheet4.xml): currently marked as the tab-selected sheet; the new Portugal sheet
# will take over that role, so just leave Swiss's own selection as-is (B5) and let adding /
# activating the new sheet clear tabSelected from Swiss automatically.

# --- Add the new "Portugal" worksheet ----------------------------------------------------------
$wsSwiss = $wb.Worksheets.Item("Swiss")
$wsPortugal = $wb.Worksheets.Add($null, $wsSwiss)
$wsPortugal.Name = "Portugal"

# Copy the layout of the Swiss sheet (headers, labels, styles, merges) as a starting point
$wsSwiss.Cells.Copy() | Out-Null
$wsPortugal.Range("A1").PasteSpecial() | Out-Null
$wsPortugal.Application.CutCopyMode = $false

# Fill in the Portugal specific values
$wsPortugal.Range("B2").Value = "Portugal Market"
$wsPortugal.Range("B4").Value = "NGC-3479/T-2459/T-2460/T2461"

# Match the column widths / selection seen in the final workbook
$wsPortugal.Columns.Item(1).ColumnWidth = 23
$wsPortugal.Columns.Item(2).ColumnWidth = 16.88671875
$wsPortugal.Columns.Item(3).ColumnWidth = 13.5546875
$wsPortugal.Columns.Item(4).ColumnWidth = 14.33203125

$wsPortugal.Range("B4").Select() | Out-Null

# Make Portugal the active / visible tab
$wsPortugal.Activate()
